$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Flip several D-column exchange amounts from their old values to 0
#    (rows that record "pvi" extra-fuel-consumption and transport/processing
#    re-use amounts that are now zeroed out).
# ---------------------------------------------------------------------------
$zeroRows = @(260,272,284,337,338,339,341,342,358,359,360,362,363,379,380,381,383,384)
foreach ($r in $zeroRows) {
    $ws.Cells.Item($r, 4).Value = 0
}

# ---------------------------------------------------------------------------
# 2) Append a brand-new "pavement structure" activity block (rows 386-398)
#    plus a handful of trailing spacer rows (399-405), mirroring the layout
#    used by the other activity blocks already on the sheet.
# ---------------------------------------------------------------------------

# -- activity header (rows 386-393) ----------------------------------------
$ws.Cells.Item(386, 1).Value = "Activity"
$ws.Cells.Item(386, 2).Value = "pavement structure"

$ws.Cells.Item(387, 1).Value = "reference product"
$ws.Cells.Item(387, 2).Value = "pavement structure"

$ws.Cells.Item(388, 1).Value = "code"
$ws.Cells.Item(388, 2).Value = "pavement_structure"

$ws.Cells.Item(389, 1).Value = "location"
$ws.Cells.Item(389, 2).Value = "NL"

$ws.Cells.Item(390, 1).Value = "amount"
$ws.Cells.Item(390, 2).Value = 1

$ws.Cells.Item(391, 1).Value = "unit"
$ws.Cells.Item(391, 2).Value = "kilometer"

$ws.Cells.Item(392, 1).Value = "original_ConversionDem2FU"
$ws.Cells.Item(392, 2).Value = 1

$ws.Cells.Item(393, 1).Value = "Exchanges"

# -- exchanges table header (row 394) ---------------------------------------
$ws.Cells.Item(394, 1).Value = "name"
$ws.Cells.Item(394, 2).Value = "reference product"
$ws.Cells.Item(394, 3).Value = "location"
$ws.Cells.Item(394, 4).Value = "amount"
$ws.Cells.Item(394, 5).Value = "unit"
$ws.Cells.Item(394, 6).Value = "database"
$ws.Cells.Item(394, 7).Value = "type"
$ws.Cells.Item(394, 8).Value = "categories"
$ws.Cells.Item(394, 9).Value = "comments"
$ws.Cells.Item(394, 10).Value = "uncertainty_type"
$ws.Cells.Item(394, 11).Value = "loc"
$ws.Cells.Item(394, 12).Value = "scale"

# -- self-referencing "production" row (row 395) -----------------------------
$ws.Cells.Item(395, 1).Formula = "=B386"
$ws.Cells.Item(395, 2).Formula = "=B387"
$ws.Cells.Item(395, 3).Value = "NL"
$ws.Cells.Item(395, 4).Value = 1
$ws.Cells.Item(395, 5).Formula = "=B391"
$ws.Cells.Item(395, 6).Value = "asphalt"
$ws.Cells.Item(395, 7).Value = "production"
$ws.Cells.Item(395, 8).Value = "(unknown)"
$ws.Cells.Item(395, 10).Value = 0

# -- AC Surf / AC Bin (x2) material exchanges (rows 396-398) -----------------
$ws.Cells.Item(396, 1).Formula = "=A332"
$ws.Cells.Item(396, 2).Formula = "=B332"
$ws.Cells.Item(396, 3).Formula = "=C332"
$ws.Cells.Item(396, 4).Formula = "=(3.75*6*0.038*1000)*2350"
$ws.Cells.Item(396, 5).Value = "kilogram"
$ws.Cells.Item(396, 6).Value = "asphalt"
$ws.Cells.Item(396, 7).Value = "technosphere"
$ws.Cells.Item(396, 8).Value = "(unknown)"
$ws.Cells.Item(396, 9).Value = "AC_surf: density 2350 kg/m3"
$ws.Cells.Item(396, 10).Value = 0

$ws.Cells.Item(397, 1).Formula = "=A353"
$ws.Cells.Item(397, 2).Formula = "=B353"
$ws.Cells.Item(397, 3).Formula = "=C353"
$ws.Cells.Item(397, 4).Formula = "=(3.75*6*0.051*1000)*2370"
$ws.Cells.Item(397, 5).Value = "kilogram"
$ws.Cells.Item(397, 6).Value = "asphalt"
$ws.Cells.Item(397, 7).Value = "technosphere"
$ws.Cells.Item(397, 8).Value = "(unknown)"
$ws.Cells.Item(397, 9).Value = "AC_bin: density 2370 kg/m3"
$ws.Cells.Item(397, 10).Value = 0

$ws.Cells.Item(398, 1).Formula = "=A353"
$ws.Cells.Item(398, 2).Formula = "=B353"
$ws.Cells.Item(398, 3).Formula = "=C353"
$ws.Cells.Item(398, 4).Formula = "=(3.75*6*0.254*1000)*2370"
$ws.Cells.Item(398, 5).Value = "kilogram"
$ws.Cells.Item(398, 6).Value = "asphalt"
$ws.Cells.Item(398, 7).Value = "technosphere"
$ws.Cells.Item(398, 8).Value = "(unknown)"
$ws.Cells.Item(398, 9).Value = "AC_base: density 2370 kg/m3"
$ws.Cells.Item(398, 10).Value = 0
$ws.Cells.Item(398, 10).Font.Size = 10
$ws.Cells.Item(398, 10).Font.Color = 0

# -- trailing spacer rows (F only, bold/black-font style) --------------------
foreach ($r in 386..388) {
    $ws.Cells.Item($r, 6).Font.Size = 10
    $ws.Cells.Item($r, 6).Font.Color = 0
}
foreach ($r in 399..405) {
    $ws.Cells.Item($r, 6).Font.Size = 10
    $ws.Cells.Item($r, 6).Font.Color = 0
}

# ---------------------------------------------------------------------------
# 3) Widen column I to fit the new "density" comments, and move the viewport
#    down to the newly added block.
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 21.5

$ws.Application.ActiveWindow.ScrollRow = 364
$ws.Range("E401").Select()
